$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The decision-table column header in cell C7 ("From") was renamed to
# "From111" when this revision was restored.
$ws.Range("C7").Value = "From111"
